$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 477.3871
$ws.Range("I80").Value = 300.84616
$ws.Range("J80").Value = 604.8889
$ws.Range("K80").Value = 902.5384799999999
$ws.Range("L80").Value = 1814.6667
$ws.Range("M80").Value = 95.46152000000006
$ws.Range("N80").Value = -3810.6667

$ws.Range("H83").Value = 477.3871
$ws.Range("I83").Value = 300.84616
$ws.Range("J83").Value = 604.8889
$ws.Range("K83").Value = 2707.61544
$ws.Range("L83").Value = 5444.0001
$ws.Range("M83").Value = 2284.38456
$ws.Range("N83").Value = -15428.0001

$ws.Range("H86").Value = 3081.6667
$ws.Range("I86").Value = 2995
$ws.Range("J86").Value = 3125
$ws.Range("K86").Value = 2995
$ws.Range("L86").Value = 3125
$ws.Range("M86").Value = -1872
$ws.Range("N86").Value = -5371

$ws.Range("H89").Value = 3081.6667
$ws.Range("I89").Value = 2995
$ws.Range("J89").Value = 3125
$ws.Range("K89").Value = 14975
$ws.Range("L89").Value = 15625
$ws.Range("M89").Value = -9359
$ws.Range("N89").Value = -26857

$ws.Range("H135").Value = 1210.7333
$ws.Range("I135").Value = 1156.5
$ws.Range("J135").Value = 1427.6666
$ws.Range("K135").Value = 10408.5
$ws.Range("L135").Value = 12848.9994
$ws.Range("M135").Value = -7873.5
$ws.Range("N135").Value = -17918.9994

$ws.Range("H138").Value = 2467.6191
$ws.Range("I138").Value = 1401.7693
$ws.Range("J138").Value = 4199.625
$ws.Range("K138").Value = 4205.3079
$ws.Range("L138").Value = 12598.875
$ws.Range("M138").Value = 934.6921000000002
$ws.Range("N138").Value = -22878.875

$ws.Range("H141").Value = 1769.25
$ws.Range("I141").Value = 1293.4286
$ws.Range("K141").Value = 3880.2858
$ws.Range("M141").Value = 1299.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2857.0557
$ws.Range("I110").Value = 2558.1428
$ws.Range("J110").Value = 3903.25
$ws.Range("K110").Value = 2558.1428
$ws.Range("L110").Value = 3903.25
$ws.Range("M110").Value = -513.1428000000001
$ws.Range("N110").Value = -7993.25

$ws.Range("H132").Value = 3002.8333
$ws.Range("I132").Value = 1576
$ws.Range("K132").Value = 4728
$ws.Range("M132").Value = -2198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2296.4167
$ws.Range("J20").Value = 1742.8572
$ws.Range("L20").Value = 1742.8572
$ws.Range("N20").Value = -2236.8572

$ws.Range("H86").Value = 3999.6667
$ws.Range("I86").Value = 2499.5
$ws.Range("K86").Value = 2499.5
$ws.Range("M86").Value = -1376.5

$ws.Range("H89").Value = 3999.6667
$ws.Range("I89").Value = 2499.5
$ws.Range("K89").Value = 12497.5
$ws.Range("M89").Value = -6881.5

$ws.Range("H134").Value = 3361.3572
$ws.Range("I134").Value = 2850.6924
$ws.Range("K134").Value = 8552.0772
$ws.Range("M134").Value = -6017.0772

$ws.Range("H140").Value = 136926.67
$ws.Range("J140").Value = 136926.67
$ws.Range("L140").Value = 136926.67
$ws.Range("N140").Value = -147286.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5076.6
$ws.Range("I31").Value = 2447.1667
$ws.Range("J31").Value = 6829.5557
$ws.Range("K31").Value = 2447.1667
$ws.Range("L31").Value = 6829.5557
$ws.Range("M31").Value = -2152.1667
$ws.Range("N31").Value = -7419.5557

$ws.Range("H34").Value = 5076.6
$ws.Range("I34").Value = 2447.1667
$ws.Range("J34").Value = 6829.5557
$ws.Range("K34").Value = 2447.1667
$ws.Range("L34").Value = 6829.5557
$ws.Range("M34").Value = -2245.1667
$ws.Range("N34").Value = -7233.5557

$ws.Range("H58").Value = 4697.8184
$ws.Range("I58").Value = 929.25
$ws.Range("K58").Value = 929.25
$ws.Range("M58").Value = -726.25

$ws.Range("H107").Value = 351.5625
$ws.Range("I107").Value = 352.72726
$ws.Range("J107").Value = 349
$ws.Range("K107").Value = 352.72726
$ws.Range("L107").Value = 349
$ws.Range("M107").Value = 1567.27274
$ws.Range("N107").Value = -4189

$ws.Range("H132").Value = 4645.5
$ws.Range("I132").Value = 3770.5
$ws.Range("K132").Value = 11311.5
$ws.Range("M132").Value = -8781.5

$ws.Range("H136").Value = 4697.8184
$ws.Range("I136").Value = 929.25
$ws.Range("K136").Value = 2787.75
$ws.Range("M136").Value = -237.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 952.3333
$ws.Range("J68").Value = 1028.5
$ws.Range("L68").Value = 3085.5
$ws.Range("N68").Value = -4707.5

$ws.Range("H71").Value = 952.3333
$ws.Range("J71").Value = 1028.5
$ws.Range("L71").Value = 9256.5
$ws.Range("N71").Value = -17368.5

$ws.Range("H94").Value = 8000
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H112").Value = 2000
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H121").Value = 1418
$ws.Range("J121").Value = 2672.1667
$ws.Range("L121").Value = 8016.500100000001
$ws.Range("N121").Value = -10636.5001

$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 3000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 9000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3900
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1566.6666
$ws.Range("I80").Value = 1100
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 1100
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -102
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 1566.6666
$ws.Range("I83").Value = 1100
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 5500
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -508
$ws.Range("N83").Value = -22484

$ws.Range("H102").Value = 1331.7222
$ws.Range("I102").Value = 1368.3529
$ws.Range("K102").Value = 1368.3529
$ws.Range("M102").Value = 253.6470999999999

$ws.Range("H122").Value = 501610.8
$ws.Range("I122").Value = 557178.7
$ws.Range("K122").Value = 1671536.1
$ws.Range("M122").Value = -1669086.1

$ws.Range("H132").Value = 50335.523
$ws.Range("I132").Value = 83754.30499999999
$ws.Range("J132").Value = 6891.1
$ws.Range("K132").Value = 251262.915
$ws.Range("L132").Value = 20673.3
$ws.Range("M132").Value = -248732.915
$ws.Range("N132").Value = -25733.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8349
$ws.Range("I7").Value = 6549.5
$ws.Range("J7").Value = 8798.875
$ws.Range("K7").Value = 6549.5
$ws.Range("L7").Value = 8798.875
$ws.Range("M7").Value = -6437.5
$ws.Range("N7").Value = -9022.875

$ws.Range("H40").Value = 5856.636
$ws.Range("I40").Value = 3669.2856
$ws.Range("J40").Value = 9684.5
$ws.Range("K40").Value = 3669.2856
$ws.Range("L40").Value = 9684.5
$ws.Range("M40").Value = -3533.2856
$ws.Range("N40").Value = -9956.5

$ws.Range("H68").Value = 8602.071
$ws.Range("I68").Value = 6487.25
$ws.Range("K68").Value = 6487.25
$ws.Range("M68").Value = -5738.25

$ws.Range("H71").Value = 8602.071
$ws.Range("I71").Value = 6487.25
$ws.Range("K71").Value = 32436.25
$ws.Range("M71").Value = -28692.25

$ws.Range("H126").Value = 8349
$ws.Range("I126").Value = 6549.5
$ws.Range("J126").Value = 8798.875
$ws.Range("K126").Value = 19648.5
$ws.Range("L126").Value = 26396.625
$ws.Range("M126").Value = -17178.5
$ws.Range("N126").Value = -31336.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 151540
$ws.Range("I4").Value = 151540
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 151540
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -151427
$ws.Range("N4").ClearContents()

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H62").Value = 11998.333
$ws.Range("J62").Value = 11998.333
$ws.Range("L62").Value = 11998.333
$ws.Range("N62").Value = -13246.333

$ws.Range("H65").Value = 11998.333
$ws.Range("J65").Value = 11998.333
$ws.Range("L65").Value = 59991.665
$ws.Range("N65").Value = -66231.66500000001

$ws.Range("H81").Value = 1495.5
$ws.Range("I81").Value = 1495.5
$ws.Range("K81").Value = 2991
$ws.Range("M81").Value = -1930

$ws.Range("H84").Value = 1495.5
$ws.Range("I84").Value = 1495.5
$ws.Range("K84").Value = 14955
$ws.Range("M84").Value = -9651

$ws.Range("H113").Value = 666.7059
$ws.Range("I113").Value = 732.25
$ws.Range("K113").Value = 2196.75
$ws.Range("M113").Value = -26.75

$ws.Range("H122").Value = 2491.6667
$ws.Range("I122").Value = 2487.5
$ws.Range("K122").Value = 7462.5
$ws.Range("M122").Value = -5012.5
